$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1407
$ws.Range("I15").Value = 1407
$ws.Range("K15").Value = 4221
$ws.Range("M15").Value = -4052
$ws.Range("H19").Value = 1099.5
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H32").Value = 1664.75
$ws.Range("I32").Value = 1599.5
$ws.Range("K32").Value = 1599.5
$ws.Range("M32").Value = -1273.5
$ws.Range("H33").Value = 161.125
$ws.Range("I33").Value = 177
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 177
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = 52
$ws.Range("N33").Value = -508
$ws.Range("H40").Value = 10834
$ws.Range("I40").Value = 7001
$ws.Range("J40").Value = 18500
$ws.Range("K40").Value = 7001
$ws.Range("L40").Value = 18500
$ws.Range("M40").Value = -6826
$ws.Range("N40").Value = -18850
$ws.Range("H137").Value = 3517.7368
$ws.Range("I137").Value = 3471
$ws.Range("K137").Value = 10413
$ws.Range("M137").Value = -7863
$ws.Range("H138").Value = 4457.0435
$ws.Range("I138").Value = 1567.12
$ws.Range("J138").Value = 5535.373
$ws.Range("K138").Value = 4701.36
$ws.Range("L138").Value = 16606.119
$ws.Range("M138").Value = 438.6400000000003
$ws.Range("N138").Value = -26886.119

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3016.0715
$ws.Range("I61").Value = 1774.6818
$ws.Range("J61").Value = 7567.8335
$ws.Range("K61").Value = 1774.6818
$ws.Range("L61").Value = 7567.8335
$ws.Range("M61").Value = -1562.6818
$ws.Range("N61").Value = -7991.8335
$ws.Range("H63").Value = 5490.6113
$ws.Range("I63").Value = 2479.125
$ws.Range("J63").Value = 7899.8
$ws.Range("K63").Value = 2479.125
$ws.Range("L63").Value = 7899.8
$ws.Range("M63").Value = -1793.125
$ws.Range("N63").Value = -9271.799999999999
$ws.Range("H66").Value = 5490.6113
$ws.Range("I66").Value = 2479.125
$ws.Range("J66").Value = 7899.8
$ws.Range("K66").Value = 12395.625
$ws.Range("L66").Value = 39499
$ws.Range("M66").Value = -8963.625
$ws.Range("N66").Value = -46363
$ws.Range("H110").Value = 140320.14
$ws.Range("I110").Value = 152987.88
$ws.Range("K110").Value = 152987.88
$ws.Range("M110").Value = -150942.88
$ws.Range("H122").Value = 5323.3
$ws.Range("I122").Value = 4002.6667
$ws.Range("K122").Value = 12008.0001
$ws.Range("M122").Value = -9558.000100000001
$ws.Range("H132").Value = 3123.64
$ws.Range("I132").Value = 3472.7334
$ws.Range("K132").Value = 10418.2002
$ws.Range("M132").Value = -7888.200199999999
$ws.Range("H136").Value = 3016.0715
$ws.Range("I136").Value = 1774.6818
$ws.Range("J136").Value = 7567.8335
$ws.Range("K136").Value = 5324.0454
$ws.Range("L136").Value = 22703.5005
$ws.Range("M136").Value = -2774.0454
$ws.Range("N136").Value = -27803.5005

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1790.4667
$ws.Range("I107").Value = 1443.3334
$ws.Range("J107").Value = 2311.1667
$ws.Range("K107").Value = 1443.3334
$ws.Range("L107").Value = 2311.1667
$ws.Range("M107").Value = 476.6666
$ws.Range("N107").Value = -6151.1667
$ws.Range("H134").Value = 28889.078
$ws.Range("I134").Value = 1920.1471
$ws.Range("K134").Value = 5760.4413
$ws.Range("M134").Value = -3225.4413

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 4974.5
$ws.Range("J13").Value = 4974.5
$ws.Range("L13").Value = 4974.5
$ws.Range("N13").Value = -5252.5
$ws.Range("H14").Value = 6890.8335
$ws.Range("J14").Value = 6890.8335
$ws.Range("L14").Value = 6890.8335
$ws.Range("N14").Value = -7230.8335
$ws.Range("H132").Value = 2562.6667
$ws.Range("I132").Value = 1444.4615
$ws.Range("K132").Value = 4333.3845
$ws.Range("M132").Value = -1803.3845
$ws.Range("H134").Value = 348159.38
$ws.Range("I134").Value = 3204.96
$ws.Range("K134").Value = 9614.880000000001
$ws.Range("M134").Value = -7079.880000000001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 833605.5600000001
$ws.Range("I8").Value = 833605.5600000001
$ws.Range("K8").Value = 2500816.68
$ws.Range("M8").Value = -2500677.68
$ws.Range("H55").Value = 44805.555
$ws.Range("I55").Value = 3250
$ws.Range("J55").Value = 50000
$ws.Range("K55").Value = 9750
$ws.Range("L55").Value = 150000
$ws.Range("M55").Value = -9573
$ws.Range("N55").Value = -150354
$ws.Range("H60").Value = 333.22223
$ws.Range("I60").Value = 271.2857
$ws.Range("K60").Value = 813.8571000000001
$ws.Range("M60").Value = -562.8571000000001
$ws.Range("H61").Value = 100
$ws.Range("I61").Value = 100
$ws.Range("K61").Value = 300
$ws.Range("M61").Value = -85
$ws.Range("H139").Value = 5849.6
$ws.Range("I139").Value = 4606.1577
$ws.Range("K139").Value = 13818.4731
$ws.Range("M139").Value = -8678.473099999999
$ws.Range("H140").Value = 2160.3333
$ws.Range("I140").Value = 1423.8276
$ws.Range("K140").Value = 4271.4828
$ws.Range("M140").Value = 908.5172000000002

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23265
$ws.Range("I15").Value = 19895
$ws.Range("J15").Value = 24950
$ws.Range("K15").Value = 19895
$ws.Range("L15").Value = 24950
$ws.Range("M15").Value = -19607
$ws.Range("N15").Value = -25526
$ws.Range("H17").Value = 3625.75
$ws.Range("J17").Value = 3625.75
$ws.Range("L17").Value = 3625.75
$ws.Range("N17").Value = -3961.75
$ws.Range("H80").Value = 1826312.4
$ws.Range("I80").Value = 1436926.2
$ws.Range("J80").Value = 2507738
$ws.Range("K80").Value = 1436926.2
$ws.Range("L80").Value = 2507738
$ws.Range("M80").Value = -1435928.2
$ws.Range("N80").Value = -2509734
$ws.Range("H81").Value = 23265
$ws.Range("I81").Value = 19895
$ws.Range("J81").Value = 24950
$ws.Range("K81").Value = 19895
$ws.Range("L81").Value = 24950
$ws.Range("M81").Value = -18897
$ws.Range("N81").Value = -26946
$ws.Range("H83").Value = 1826312.4
$ws.Range("I83").Value = 1436926.2
$ws.Range("J83").Value = 2507738
$ws.Range("K83").Value = 7184631
$ws.Range("L83").Value = 12538690
$ws.Range("M83").Value = -7179639
$ws.Range("N83").Value = -12548674
$ws.Range("H84").Value = 23265
$ws.Range("I84").Value = 19895
$ws.Range("J84").Value = 24950
$ws.Range("K84").Value = 59685
$ws.Range("L84").Value = 74850
$ws.Range("M84").Value = -54693
$ws.Range("N84").Value = -84834
$ws.Range("H98").Value = 49500
$ws.Range("J98").Value = 49500
$ws.Range("L98").Value = 49500
$ws.Range("N98").Value = -55490
$ws.Range("H122").Value = 5094.5
$ws.Range("I122").Value = 1996
$ws.Range("J122").Value = 6127.3335
$ws.Range("K122").Value = 5988
$ws.Range("L122").Value = 18382.0005
$ws.Range("M122").Value = -3538
$ws.Range("N122").Value = -23282.0005
$ws.Range("H126").Value = 4029.6
$ws.Range("I126").Value = 3099.6667
$ws.Range("K126").Value = 9299.000100000001
$ws.Range("M126").Value = -6829.000100000001
$ws.Range("H132").Value = 74509.87
$ws.Range("I132").Value = 9772.333000000001
$ws.Range("K132").Value = 29316.999
$ws.Range("M132").Value = -26786.999
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1179.4
$ws.Range("I16").Value = 1249.5
$ws.Range("J16").Value = 899
$ws.Range("K16").Value = 1249.5
$ws.Range("L16").Value = 899
$ws.Range("M16").Value = -1079.5
$ws.Range("N16").Value = -1239
$ws.Range("H55").Value = 739.2857
$ws.Range("I55").Value = 203.95653
$ws.Range("J55").Value = 3201.8
$ws.Range("K55").Value = 203.95653
$ws.Range("L55").Value = 3201.8
$ws.Range("M55").Value = -30.95652999999999
$ws.Range("N55").Value = -3547.8
$ws.Range("H132").Value = 3966.5
$ws.Range("I132").Value = 2266.6667
$ws.Range("K132").Value = 6800.000100000001
$ws.Range("M132").Value = -4270.000100000001
$ws.Range("H136").Value = 1825616.6
$ws.Range("I136").Value = 2506047
$ws.Range("K136").Value = 7518141
$ws.Range("M136").Value = -7515591

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 98000.27
$ws.Range("I62").Value = 257000.75
$ws.Range("J62").Value = 7142.857
$ws.Range("K62").Value = 257000.75
$ws.Range("L62").Value = 7142.857
$ws.Range("M62").Value = -256376.75
$ws.Range("N62").Value = -8390.857
$ws.Range("H65").Value = 98000.27
$ws.Range("I65").Value = 257000.75
$ws.Range("J65").Value = 7142.857
$ws.Range("K65").Value = 1285003.75
$ws.Range("L65").Value = 35714.285
$ws.Range("M65").Value = -1281883.75
$ws.Range("N65").Value = -41954.285
$ws.Range("H100").Value = 1036.9375
$ws.Range("I100").Value = 1138.1538
$ws.Range("K100").Value = 2276.3076
$ws.Range("M100").Value = -1735.3076
$ws.Range("H122").Value = 37040190
$ws.Range("I122").Value = 43480920
$ws.Range("J122").Value = 6001.25
$ws.Range("K122").Value = 130442760
$ws.Range("L122").Value = 18003.75
$ws.Range("M122").Value = -130440310
$ws.Range("N122").Value = -22903.75
$ws.Range("H125").Value = 72471.5
$ws.Range("J125").Value = 72471.5
$ws.Range("L125").Value = 72471.5
$ws.Range("N125").Value = -82311.5
$ws.Range("H132").Value = 18290.562
$ws.Range("I132").Value = 2559.65
$ws.Range("J132").Value = 44508.75
$ws.Range("K132").Value = 7678.950000000001
$ws.Range("L132").Value = 133526.25
$ws.Range("M132").Value = -5148.950000000001
$ws.Range("N132").Value = -138586.25
$ws.Range("H136").Value = 9090484
$ws.Range("I136").Value = 14930385
$ws.Range("J136").Value = 135969.53
$ws.Range("K136").Value = 44791155
$ws.Range("L136").Value = 407908.59
$ws.Range("M136").Value = -44788605
$ws.Range("N136").Value = -413008.59
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
